$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "Graphical interface" block & hyperlinked bug-report row ---
$ws.Hyperlinks.Delete()
$ws.Range("B24:D25").UnMerge()
$ws.Rows("20:25").Delete()

# --- Make room for a new "invalid email" profile test case (new row 18) ---
$ws.Rows("18:18").Insert()
$ws.Range("B17:D17").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)
$ws.Rows("18").RowHeight = 16.8
$ws.Application.CutCopyMode = $false

# --- Update the "change profile" test case texts ---
$ws.Range("B15").Value = "Именить User name : Alex_323 ; Phone number : +7 4954434343"
$ws.Range("B18").Value = "Именить email address: testcaseqa.com; Phone number : +7 4954434343"

# --- Add the new "invalid phone" profile test case in the row freed up by the earlier delete ---
$ws.Range("B20:D20").Copy()
$ws.Range("B21:D21").PasteSpecial(-4122)
$ws.Rows("21").RowHeight = 16.8
$ws.Application.CutCopyMode = $false
$ws.Range("B21").Value = "Именить Phone number на не валидный: +7 3343524334332"

# --- Fix up the border on the first "Authorization" data row so it matches the others ---
$ws.Range("C10:D10").Copy()
$ws.Range("C9:D9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Fill in the pass/fail test results ---
$ws.Range("C9").Value = "pass"
$ws.Range("D9").Value = "pass"
$ws.Range("C10").Value = "fail bug#1"
$ws.Range("D10").Value = "fail bug#1"
$ws.Range("C11").Value = "pass"
$ws.Range("D11").Value = "pass"
$ws.Range("C12").Value = "pass"
$ws.Range("D12").Value = "pass"
$ws.Range("C13").Value = "pass"
$ws.Range("D13").Value = "pass"
$ws.Range("C15").Value = "pass"
$ws.Range("D15").Value = "pass"
$ws.Range("C16").Value = "pass"
$ws.Range("D16").Value = "pass"
$ws.Range("C17").Value = "pass"
$ws.Range("D17").Value = "pass"
$ws.Range("C18").Value = "pass"
$ws.Range("D18").Value = "pass"
$ws.Range("C19").Value = "pass"
$ws.Range("D19").Value = "pass"
$ws.Range("C20").Value = "pass"
$ws.Range("D20").Value = "pass"
$ws.Range("C21").Value = "pass"
$ws.Range("D21").Value = "pass"

# --- Restore the cursor/selection shown in the workbook ---
$ws.Range("G5").Select()
